$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") for rows 2-10 changes from serial date 45208 to 45212
for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 3).Value = 45212
}
